# Auto-generated script: apply 2023-11-16 data update to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6766
$ws.Range("J3").Value = 7149
$ws.Range("G4").Value = 1472
$ws.Range("I4").Value = 1774
$ws.Range("J4").Value = 1557
$ws.Range("J5").Value = 560
$ws.Range("J6").Value = 9531
$ws.Range("G7").Value = 24697
$ws.Range("I7").Value = 26232
$ws.Range("J7").Value = 25563

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 115
$ws.Range("J7").Value = 741
$ws.Range("J8").Value = 1605
$ws.Range("J10").Value = 186
$ws.Range("J11").Value = 440
$ws.Range("J12").Value = 54
$ws.Range("J18").Value = 215
$ws.Range("J19").Value = 747
$ws.Range("J30").Value = 91
$ws.Range("J31").Value = 255
$ws.Range("J33").Value = 1154
$ws.Range("J36").Value = 348
$ws.Range("J37").Value = 791
$ws.Range("J39").Value = 17
$ws.Range("J41").Value = 178
$ws.Range("J42").Value = 1102
$ws.Range("J43").Value = 219
$ws.Range("J44").Value = 195
$ws.Range("J52").Value = 646
$ws.Range("J53").Value = 366
$ws.Range("J54").Value = 491
$ws.Range("J55").Value = 393
$ws.Range("J60").Value = 148
$ws.Range("G63").Value = 274
$ws.Range("J63").Value = 81
$ws.Range("J65").Value = 639
$ws.Range("J67").Value = 959
$ws.Range("J76").Value = 376
$ws.Range("J78").Value = 300
$ws.Range("J79").Value = 720
$ws.Range("J83").Value = 511
$ws.Range("J84").Value = 214
$ws.Range("I86").Value = 169
$ws.Range("J89").Value = 327
$ws.Range("J90").Value = 273
$ws.Range("J91").Value = 294
$ws.Range("J92").Value = 80
$ws.Range("J94").Value = 272
$ws.Range("J95").Value = 370
$ws.Range("J99").Value = 392
$ws.Range("J100").Value = 47
$ws.Range("G101").Value = 24697
$ws.Range("I101").Value = 26232
$ws.Range("J101").Value = 25563

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 226
$ws.Range("J7").Value = 741

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 201
$ws.Range("J7").Value = 440

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 99
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 327

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 153
$ws.Range("J6").Value = 275
$ws.Range("J7").Value = 646

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 366

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 434
$ws.Range("J3").Value = 481
$ws.Range("J6").Value = 564
$ws.Range("J7").Value = 1605

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 151
$ws.Range("J7").Value = 511

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 268
$ws.Range("J3").Value = 384
$ws.Range("J6").Value = 401
$ws.Range("J7").Value = 1154

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 127
$ws.Range("J3").Value = 132
$ws.Range("J7").Value = 370

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 233
$ws.Range("J4").Value = 31
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 791

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 187
$ws.Range("J3").Value = 181
$ws.Range("J7").Value = 639

# Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 156
$ws.Range("J7").Value = 392

# Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 91

# Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 89
$ws.Range("J7").Value = 255

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 355
$ws.Range("J6").Value = 268
$ws.Range("J7").Value = 959

# South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 214

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 102
$ws.Range("J4").Value = 36
$ws.Range("J7").Value = 491

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 181
$ws.Range("J6").Value = 289
$ws.Range("J7").Value = 747

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 195

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 79
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 376

# Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 178

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 217
$ws.Range("J7").Value = 1102

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 105
$ws.Range("J7").Value = 186

# Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 94
$ws.Range("J7").Value = 300

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 76
$ws.Range("J6").Value = 219
$ws.Range("J7").Value = 393

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 294

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 200
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 720

# Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 43
$ws.Range("J7").Value = 215

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 112
$ws.Range("J7").Value = 348

# Wrigleyville
$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 47

# West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 272

# Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 25
$ws.Range("J6").Value = 118

# Greektown
$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("J2").Value = 4
$ws.Range("J6").Value = 17

# West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 80

# Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 80
$ws.Range("I7").Value = 169

# Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 96
$ws.Range("J7").Value = 273

# Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 52
$ws.Range("J7").Value = 148

# Hyde Park
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 130
$ws.Range("J7").Value = 219

# Archer Heights
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J2").Value = 36
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 115

# Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 54
